# Update countries & provincias Spain
#
# Daily refresh of the COVID-19 "paises" (countries) table on sheet "Pais".
# The "last updated" banner in A1 moves from 17:52 to 18:22, most countries'
# stat columns (B:H = Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) get fresh totals, and a few rows swap
# places because the refreshed "Casos totales" numbers changed their
# descending sort order (Chile/Dinamarca, Ecuador/India, Argelia/Singapur/
# Ucrania, Cuba/Burkina Faso/Uzbekistan/Oman).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 18:22"

# --- Country rows (A:H) --------------------------------------------------
# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Range("B4").Value = 323516
$ws.Range("C4").Value = 12159
$ws.Range("D4").Value = 16590
$ws.Range("E4").Value = 297755
$ws.Range("F4").Value = 8474
$ws.Range("G4").Value = 719
$ws.Range("H4").Value = 9171

# Row 6: Italia
$ws.Cells.Item(6, 1).Value = "Italia"
$ws.Range("B6").Value = 128948
$ws.Range("C6").Value = 4316
$ws.Range("D6").Value = 21815
$ws.Range("E6").Value = 91246
$ws.Range("F6").Value = 3977
$ws.Range("G6").Value = 525
$ws.Range("H6").Value = 15887

# Row 7: Alemania
$ws.Cells.Item(7, 1).Value = "Alemania"
$ws.Range("B7").Value = 98737
$ws.Range("C7").Value = 2645
$ws.Range("D7").Value = 26400
$ws.Range("E7").Value = 70814
$ws.Range("F7").Value = 3936
$ws.Range("G7").Value = 79
$ws.Range("H7").Value = 1523

# Row 17: Austria
$ws.Cells.Item(17, 1).Value = "Austria"
$ws.Range("B17").Value = 11934
$ws.Range("C17").Value = 153
$ws.Range("D17").Value = 2998
$ws.Range("E17").Value = 8732
$ws.Range("F17").Value = 244
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 204

# Row 24: Noruega
$ws.Cells.Item(24, 1).Value = "Noruega"
$ws.Range("B24").Value = 5686
$ws.Range("C24").Value = 136
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 5584
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 8
$ws.Range("H24").Value = 70

# Row 27: Chequia
$ws.Cells.Item(27, 1).Value = "Chequia"
$ws.Range("B27").Value = 4543
$ws.Range("C27").Value = 71
$ws.Range("D27").Value = 96
$ws.Range("E27").Value = 4380
$ws.Range("F27").Value = 86
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 67

# Row 28: Chile (was Dinamarca's row; Chile overtakes it in rank)
$ws.Cells.Item(28, 1).Value = "Chile"
$ws.Range("B28").Value = 4471
$ws.Range("C28").Value = 310
$ws.Range("D28").Value = 528
$ws.Range("E28").Value = 3909
$ws.Range("F28").Value = 38
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 34

# Row 29: Dinamarca (was Chile's row)
$ws.Cells.Item(29, 1).Value = "Dinamarca"
$ws.Range("B29").Value = 4369
$ws.Range("C29").Value = 292
$ws.Range("D29").Value = 1327
$ws.Range("E29").Value = 2863
$ws.Range("F29").Value = 142
$ws.Range("G29").Value = 18
$ws.Range("H29").Value = 179

# Row 33: Ecuador (was India's row; Ecuador overtakes it in rank)
$ws.Cells.Item(33, 1).Value = "Ecuador"
$ws.Range("B33").Value = 3646
$ws.Range("C33").Value = 181
$ws.Range("D33").Value = 100
$ws.Range("E33").Value = 3366
$ws.Range("F33").Value = 100
$ws.Range("G33").Value = 8
$ws.Range("H33").Value = 180

# Row 34: India (was Ecuador's row)
$ws.Cells.Item(34, 1).Value = "India"
$ws.Range("B34").Value = 3588
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 229
$ws.Range("E34").Value = 3260
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 99

# Row 55: Argelia (was Singapur's row; Argelia overtakes it in rank)
$ws.Cells.Item(55, 1).Value = "Argelia"
$ws.Range("B55").Value = 1320
$ws.Range("C55").Value = 69
$ws.Range("D55").Value = 90
$ws.Range("E55").Value = 1078
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 22
$ws.Range("H55").Value = 152

# Row 56: Singapur (was Ucrania's row)
$ws.Cells.Item(56, 1).Value = "Singapur"
$ws.Range("B56").Value = 1309
$ws.Range("C56").Value = 120
$ws.Range("D56").Value = 297
$ws.Range("E56").Value = 1006
$ws.Range("F56").Value = 24
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 6

# Row 57: Ucrania (was Argelia's row)
$ws.Cells.Item(57, 1).Value = "Ucrania"
$ws.Range("B57").Value = 1251
$ws.Range("C57").Value = 26
$ws.Range("D57").Value = 25
$ws.Range("E57").Value = 1194
$ws.Range("F57").Value = 16
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 32

# Row 93: Cuba (was Burkina Faso's row; Cuba overtakes it in rank)
$ws.Cells.Item(93, 1).Value = "Cuba"
$ws.Range("B93").Value = 320
$ws.Range("C93").Value = 32
$ws.Range("D93").Value = 15
$ws.Range("E93").Value = 297
$ws.Range("F93").Value = 11
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 8

# Row 94: Burkina Faso (was Uzbekistan's row)
$ws.Cells.Item(94, 1).Value = "Burkina Faso"
$ws.Range("B94").Value = 318
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 66
$ws.Range("E94").Value = 236
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 16

# Row 95: Uzbekistan (was Oman's row)
$ws.Cells.Item(95, 1).Value = "Uzbekistan"
$ws.Range("B95").Value = 298
$ws.Range("C95").Value = 32
$ws.Range("D95").Value = 30
$ws.Range("E95").Value = 266
$ws.Range("F95").Value = 8
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 2

# Row 96: Oman (was Cuba's row)
$ws.Cells.Item(96, 1).Value = "Oman"
$ws.Range("B96").Value = 298
$ws.Range("C96").Value = 21
$ws.Range("D96").Value = 61
$ws.Range("E96").Value = 235
$ws.Range("F96").Value = 3
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 2
